$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 80.181816
$ws.Range("I11").Value = 80.181816
$ws.Range("K11").Value = 80.181816
$ws.Range("M11").Value = 59.818184
$ws.Range("H17").Value = 1339.0588
$ws.Range("J17").Value = 1339.0588
$ws.Range("L17").Value = 4017.1764
$ws.Range("N17").Value = -4353.1764
$ws.Range("H64").Value = 57372.438
$ws.Range("I64").Value = 85247.375
$ws.Range("J64").Value = 29497.5
$ws.Range("K64").Value = 85247.375
$ws.Range("L64").Value = 29497.5
$ws.Range("M64").Value = -84999.375
$ws.Range("N64").Value = -29993.5
$ws.Range("H67").Value = 57372.438
$ws.Range("I67").Value = 85247.375
$ws.Range("J67").Value = 29497.5
$ws.Range("K67").Value = 85247.375
$ws.Range("L67").Value = 29497.5
$ws.Range("M67").Value = -84389.375
$ws.Range("N67").Value = -31213.5
$ws.Range("H132").Value = 5138.073
$ws.Range("I132").Value = 5091.5713
$ws.Range("J132").Value = 5409.3335
$ws.Range("K132").Value = 15274.7139
$ws.Range("L132").Value = 16228.0005
$ws.Range("M132").Value = -12744.7139
$ws.Range("N132").Value = -21288.0005
$ws.Range("H138").Value = 242187.95
$ws.Range("I138").Value = 628195.5600000001
$ws.Range("K138").Value = 1884586.68
$ws.Range("M138").Value = -1879446.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 123934.06
$ws.Range("I45").Value = 158989.61
$ws.Range("K45").Value = 158989.61
$ws.Range("M45").Value = -158612.61
$ws.Range("H63").Value = 1010.7143
$ws.Range("I63").Value = 1029.3334
$ws.Range("K63").Value = 1029.3334
$ws.Range("M63").Value = -343.3334
$ws.Range("H66").Value = 1010.7143
$ws.Range("I66").Value = 1029.3334
$ws.Range("K66").Value = 5146.666999999999
$ws.Range("M66").Value = -1714.666999999999
$ws.Range("H132").Value = 3129.712
$ws.Range("I132").Value = 2191.366
$ws.Range("K132").Value = 6574.098
$ws.Range("M132").Value = -4044.098

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6835265
$ws.Range("I99").Value = 11615401
$ws.Range("J99").Value = 6499.143
$ws.Range("K99").Value = 11615401
$ws.Range("L99").Value = 6499.143
$ws.Range("M99").Value = -11613903
$ws.Range("N99").Value = -9495.143
$ws.Range("H122").Value = 21576
$ws.Range("I122").Value = 120000
$ws.Range("J122").Value = 1891.2
$ws.Range("K122").Value = 360000
$ws.Range("L122").Value = 5673.6
$ws.Range("M122").Value = -357550
$ws.Range("N122").Value = -10573.6
$ws.Range("H126").Value = 6835265
$ws.Range("I126").Value = 11615401
$ws.Range("J126").Value = 6499.143
$ws.Range("K126").Value = 34846203
$ws.Range("L126").Value = 19497.429
$ws.Range("M126").Value = -34843733
$ws.Range("N126").Value = -24437.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 11170.615
$ws.Range("I68").Value = 2625
$ws.Range("J68").Value = 12724.363
$ws.Range("K68").Value = 7875
$ws.Range("L68").Value = 38173.089
$ws.Range("M68").Value = -7064
$ws.Range("N68").Value = -39795.089
$ws.Range("H71").Value = 11170.615
$ws.Range("I71").Value = 2625
$ws.Range("J71").Value = 12724.363
$ws.Range("K71").Value = 23625
$ws.Range("L71").Value = 114519.267
$ws.Range("M71").Value = -19569
$ws.Range("N71").Value = -122631.267
$ws.Range("H80").Value = 86845.5
$ws.Range("J80").Value = 86845.5
$ws.Range("L80").Value = 260536.5
$ws.Range("N80").Value = -262408.5
$ws.Range("H83").Value = 86845.5
$ws.Range("J83").Value = 86845.5
$ws.Range("L83").Value = 781609.5
$ws.Range("N83").Value = -790969.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6679.276
$ws.Range("I80").Value = 6363.96
$ws.Range("J80").Value = 8650
$ws.Range("K80").Value = 6363.96
$ws.Range("L80").Value = 8650
$ws.Range("M80").Value = -5365.96
$ws.Range("N80").Value = -10646
$ws.Range("H83").Value = 6679.276
$ws.Range("I83").Value = 6363.96
$ws.Range("J83").Value = 8650
$ws.Range("K83").Value = 31819.8
$ws.Range("L83").Value = 43250
$ws.Range("M83").Value = -26827.8
$ws.Range("N83").Value = -53234
$ws.Range("H97").Value = 6206.2607
$ws.Range("I97").Value = 7598.6665
$ws.Range("K97").Value = 7598.6665
$ws.Range("M97").Value = -7102.6665
$ws.Range("H122").Value = 8190.3
$ws.Range("I122").Value = 5277.636
$ws.Range("K122").Value = 15832.908
$ws.Range("M122").Value = -13382.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11849.842
$ws.Range("J22").Value = 1835.5454
$ws.Range("L22").Value = 1835.5454
$ws.Range("N22").Value = -2425.5454
$ws.Range("H27").Value = 11849.842
$ws.Range("J27").Value = 1835.5454
$ws.Range("L27").Value = 1835.5454
$ws.Range("N27").Value = -2049.5454
$ws.Range("H40").Value = 34687.176
$ws.Range("I40").Value = 54837.875
$ws.Range("K40").Value = 54837.875
$ws.Range("M40").Value = -54701.875
$ws.Range("H68").Value = 3056.5833
$ws.Range("I68").Value = 1888.75
$ws.Range("J68").Value = 4224.4165
$ws.Range("K68").Value = 1888.75
$ws.Range("L68").Value = 4224.4165
$ws.Range("M68").Value = -1139.75
$ws.Range("N68").Value = -5722.4165
$ws.Range("H71").Value = 3056.5833
$ws.Range("I71").Value = 1888.75
$ws.Range("J71").Value = 4224.4165
$ws.Range("K71").Value = 9443.75
$ws.Range("L71").Value = 21122.0825
$ws.Range("M71").Value = -5699.75
$ws.Range("N71").Value = -28610.0825
$ws.Range("H122").Value = 6438.1763
$ws.Range("I122").Value = 7168.625
$ws.Range("K122").Value = 21505.875
$ws.Range("M122").Value = -19055.875
$ws.Range("H136").Value = 5632.6665
$ws.Range("I136").Value = 3299.5
$ws.Range("K136").Value = 9898.5
$ws.Range("M136").Value = -7348.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 44666.332
$ws.Range("J112").Value = 44666.332
$ws.Range("L112").Value = 44666.332
$ws.Range("N112").Value = -47620.332
$ws.Range("H122").Value = 4698.7407
$ws.Range("I122").Value = 2069.4285
$ws.Range("K122").Value = 6208.2855
$ws.Range("M122").Value = -3758.2855
$ws.Range("H132").Value = 14281.625
$ws.Range("I132").Value = 17028.486
$ws.Range("J132").Value = 5042.1816
$ws.Range("K132").Value = 51085.458
$ws.Range("L132").Value = 15126.5448
$ws.Range("M132").Value = -48555.458
$ws.Range("N132").Value = -20186.5448
$ws.Range("H136").Value = 816097.3
$ws.Range("I136").Value = 1288999.4
$ws.Range("J136").Value = 5408
$ws.Range("K136").Value = 3866998.2
$ws.Range("L136").Value = 16224
$ws.Range("M136").Value = -3864448.2
$ws.Range("N136").Value = -21324
